# Weekly cryptos list refresh (GitHub Actions scheduled update).
# For each affected row the "Price" (column D) and "Volume(1h)" (column E)
# figures are refreshed with the latest scraped values. A handful of rows
# (18, 21, 22) also have their coin Name/Link (columns B/C) updated because
# the underlying ranking re-ordered three coins (ShibaInu /
# WrappedliquidstakedEther2.0 / Dai) among themselves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Many of the "Price" values look like plain numbers (e.g. "0.9986"); if we
# simply assign them Excel will happily reinterpret the string as a numeric
# value (and silently normalise things like trailing zeros, e.g.
# "0.04870" -> 0.0487). The source data must stay literal text, so for every
# cell we touch we temporarily force a text number-format before writing the
# value. That also stamps the cell with a new style, so afterwards we copy
# back the plain/default style taken from an unrelated, untouched cell.
$refStyle = $ws.Range("B2").Style

$updates = @(
    @{ Cell = "D2"; Value = "31.314.34" },
    @{ Cell = "E2"; Value = "  +1.77%  " },
    @{ Cell = "D3"; Value = "1.950.37" },
    @{ Cell = "E3"; Value = "  +1.15%  " },
    @{ Cell = "D4"; Value = "0.9986" },
    @{ Cell = "E4"; Value = "  -0.25%  " },
    @{ Cell = "D5"; Value = "243.29" },
    @{ Cell = "E5"; Value = "  +0.44%  " },
    @{ Cell = "D6"; Value = "0.9979" },
    @{ Cell = "E6"; Value = "  -0.23%  " },
    @{ Cell = "D7"; Value = "0.4817" },
    @{ Cell = "E7"; Value = "  -1.01%  " },
    @{ Cell = "D8"; Value = "0.2927" },
    @{ Cell = "E8"; Value = "  -0.35%  " },
    @{ Cell = "D9"; Value = "0.06834" },
    @{ Cell = "E9"; Value = "  -0.16%  " },
    @{ Cell = "D10"; Value = "20.06" },
    @{ Cell = "E10"; Value = "  +4.52%  " },
    @{ Cell = "D11"; Value = "105.02" },
    @{ Cell = "E11"; Value = "  -0.62%  " },
    @{ Cell = "D12"; Value = "1.958.46" },
    @{ Cell = "E12"; Value = "  +1.57%  " },
    @{ Cell = "D13"; Value = "0.07834" },
    @{ Cell = "E13"; Value = "  +0.89%  " },
    @{ Cell = "D14"; Value = "5.343" },
    @{ Cell = "E14"; Value = "  +0.37%  " },
    @{ Cell = "D15"; Value = "0.6922" },
    @{ Cell = "E15"; Value = "  -0.66%  " },
    @{ Cell = "D16"; Value = "300.35" },
    @{ Cell = "E16"; Value = "  +9.48%  " },
    @{ Cell = "D17"; Value = "31.335.23" },
    @{ Cell = "E17"; Value = "  +1.90%  " },
    @{ Cell = "B18"; Value = "WrappedliquidstakedEther2.0" },
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth" },
    @{ Cell = "D18"; Value = "2.219.92" },
    @{ Cell = "E18"; Value = "  +1.48%  " },
    @{ Cell = "D19"; Value = "13.05" },
    @{ Cell = "E19"; Value = "  +0.47%  " },
    @{ Cell = "D20"; Value = "5.613" },
    @{ Cell = "E20"; Value = "  +0.56%  " },
    @{ Cell = "B21"; Value = "ShibaInu" },
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib" },
    @{ Cell = "D21"; Value = "0.000007647" },
    @{ Cell = "E21"; Value = "  -0.37%  " },
    @{ Cell = "B22"; Value = "Dai" },
    @{ Cell = "C22"; Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai" },
    @{ Cell = "D22"; Value = "0.9988" },
    @{ Cell = "E22"; Value = "  -0.13%  " },
    @{ Cell = "D23"; Value = "0.9987" },
    @{ Cell = "E23"; Value = "  -0.21%  " },
    @{ Cell = "D24"; Value = "6.489" },
    @{ Cell = "E24"; Value = "  +0.50%  " },
    @{ Cell = "D25"; Value = "9.607" },
    @{ Cell = "E25"; Value = "  -2.47%  " },
    @{ Cell = "D26"; Value = "168.88" },
    @{ Cell = "E26"; Value = "  +2.66%  " },
    @{ Cell = "D27"; Value = "20.01" },
    @{ Cell = "E27"; Value = "  +3.05%  " },
    @{ Cell = "D28"; Value = "2.139" },
    @{ Cell = "E28"; Value = "  -0.53%  " },
    @{ Cell = "D29"; Value = "1.401" },
    @{ Cell = "E29"; Value = "  +1.40%  " },
    @{ Cell = "D30"; Value = "0.1019" },
    @{ Cell = "E30"; Value = "  -1.45%  " },
    @{ Cell = "D31"; Value = "4.660" },
    @{ Cell = "E31"; Value = "  +1.90%  " },
    @{ Cell = "D32"; Value = "1.539" },
    @{ Cell = "E32"; Value = "  -0.34%  " },
    @{ Cell = "D33"; Value = "4.383" },
    @{ Cell = "E33"; Value = "  +0.85%  " },
    @{ Cell = "D34"; Value = "0.04870" },
    @{ Cell = "E34"; Value = "  -0.15%  " },
    @{ Cell = "D35"; Value = "0.7481" },
    @{ Cell = "E35"; Value = "  -1.04%  " },
    @{ Cell = "D36"; Value = "1.137" },
    @{ Cell = "E36"; Value = "  -0.55%  " },
    @{ Cell = "D37"; Value = "2.722" },
    @{ Cell = "E37"; Value = "  +0.33%  " },
    @{ Cell = "D38"; Value = "0.01968" },
    @{ Cell = "E38"; Value = "  -1.08%  " },
    @{ Cell = "D39"; Value = "6.629" },
    @{ Cell = "E39"; Value = "  +2.24%  " },
    @{ Cell = "D40"; Value = "2.650" },
    @{ Cell = "E40"; Value = "  -0.15%  " },
    @{ Cell = "D41"; Value = "76.98" },
    @{ Cell = "E41"; Value = "  -1.65%  " },
    @{ Cell = "D42"; Value = "2.059" },
    @{ Cell = "E42"; Value = "  -0.23%  " },
    @{ Cell = "D43"; Value = "0.8762" },
    @{ Cell = "E43"; Value = "  -0.90%  " },
    @{ Cell = "D44"; Value = "0.4399" },
    @{ Cell = "E44"; Value = "  -0.71%  " },
    @{ Cell = "D45"; Value = "106.79" },
    @{ Cell = "E45"; Value = "  -0.86%  " },
    @{ Cell = "D46"; Value = "0.9978" },
    @{ Cell = "E46"; Value = "  -0.30%  " },
    @{ Cell = "D47"; Value = "1.014.67" },
    @{ Cell = "E47"; Value = "  +3.89%  " },
    @{ Cell = "D48"; Value = "7.630" },
    @{ Cell = "E48"; Value = "  -3.40%  " },
    @{ Cell = "D49"; Value = "9.173" },
    @{ Cell = "E49"; Value = "  -1.01%  " },
    @{ Cell = "D50"; Value = "0.1219" },
    @{ Cell = "E50"; Value = "  -1.69%  " },
    @{ Cell = "D51"; Value = "35.38" },
    @{ Cell = "E51"; Value = "  -1.99%  " }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $update.Value
    $cell.Style = $refStyle
}
